$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for the rows that changed, per repull of data.
$updates = @{
    "F2"  = -6
    "F3"  = -4
    "F4"  = -1
    "F8"  = 1
    "F13" = 1
    "F16" = -1
    "F17" = -7
    "F19" = -7
    "F20" = -3
    "F21" = -5
    "F26" = 3
    "F29" = -2
    "F30" = -4
    "F31" = 5
    "F32" = 5
    "F34" = -3
    "F35" = 1
    "F36" = 0
    "F38" = -5
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
